$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force D2:E51 to remain text (matches original inlineStr/text cell type)
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "27.524.10"
$ws.Range("E2").Value = "  +2.26%  "

$ws.Range("D3").Value = "1.847.50"
$ws.Range("E3").Value = "  +1.90%  "

$ws.Range("D4").Value = "1.012"
$ws.Range("E4").Value = "  +0.94%  "

$ws.Range("D5").Value = "314.99"
$ws.Range("E5").Value = "  +1.92%  "

$ws.Range("D6").Value = "1.012"
$ws.Range("E6").Value = "  +0.94%  "

$ws.Range("D7").Value = "0.4758"
$ws.Range("E7").Value = "  +2.21%  "

$ws.Range("D8").Value = "0.3712"
$ws.Range("E8").Value = "  +1.13%  "

$ws.Range("D9").Value = "0.07489"
$ws.Range("E9").Value = "  +1.99%  "

$ws.Range("D10").Value = "0.8899"
$ws.Range("E10").Value = "  +2.44%  "

$ws.Range("D11").Value = "20.56"
$ws.Range("E11").Value = "  +1.30%  "

$ws.Range("D12").Value = "1.848.18"
$ws.Range("E12").Value = "  +1.82%  "

$ws.Range("D13").Value = "0.07388"
$ws.Range("E13").Value = "  +4.10%  "

$ws.Range("D14").Value = "5.485"
$ws.Range("E14").Value = "  +2.41%  "

$ws.Range("D15").Value = "93.72"
$ws.Range("E15").Value = "  +2.64%  "

$ws.Range("D16").Value = "6.616"
$ws.Range("E16").Value = "  +1.76%  "

$ws.Range("D17").Value = "1.014"
$ws.Range("E17").Value = "  +1.03%  "

$ws.Range("D18").Value = "0.000008871"
$ws.Range("E18").Value = "  +2.19%  "

$ws.Range("D19").Value = "1.013"
$ws.Range("E19").Value = "  +1.06%  "

$ws.Range("D20").Value = "14.87"
$ws.Range("E20").Value = "  +1.69%  "

$ws.Range("D21").Value = "27.458.37"
$ws.Range("E21").Value = "  +1.94%  "

$ws.Range("D22").Value = "5.349"
$ws.Range("E22").Value = "  +1.02%  "

$ws.Range("D23").Value = "10.73"
$ws.Range("E23").Value = "  +1.29%  "

$ws.Range("D24").Value = "2.085.05"
$ws.Range("E24").Value = "  +1.93%  "

$ws.Range("D25").Value = "1.899"
$ws.Range("E25").Value = "  +0.42%  "

$ws.Range("D26").Value = "152.28"
$ws.Range("E26").Value = "  +1.19%  "

$ws.Range("D27").Value = "18.67"
$ws.Range("E27").Value = "  +1.93%  "

$ws.Range("D28").Value = "2.182"
$ws.Range("E28").Value = "  +0.85%  "

$ws.Range("D29").Value = "5.290"
$ws.Range("E29").Value = "  +0.25%  "

$ws.Range("D30").Value = "118.61"
$ws.Range("E30").Value = "  +2.71%  "

$ws.Range("D31").Value = "0.09007"
$ws.Range("E31").Value = "  +0.74%  "

$ws.Range("D32").Value = "0.7632"
$ws.Range("E32").Value = "  +0.71%  "

$ws.Range("D33").Value = "1.186"
$ws.Range("E33").Value = "  +2.57%  "

$ws.Range("D34").Value = "4.580"
$ws.Range("E34").Value = "  +2.24%  "

$ws.Range("D35").Value = "2.951"
$ws.Range("E35").Value = "  +0.88%  "

$ws.Range("D36").Value = "1.013"
$ws.Range("E36").Value = "  +1.02%  "

$ws.Range("D37").Value = "1.107"
$ws.Range("E37").Value = "  +2.00%  "

$ws.Range("D38").Value = "0.05361"
$ws.Range("E38").Value = "  +1.66%  "

$ws.Range("D39").Value = "0.01965"
$ws.Range("E39").Value = "  +0.58%  "

$ws.Range("D40").Value = "3.003"
$ws.Range("E40").Value = "  +0.95%  "

$ws.Range("D41").Value = "7.362"
$ws.Range("E41").Value = "  +2.54%  "

$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").Value = "2.393"
$ws.Range("E42").Value = "  +3.76%  "

$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").Value = "0.5353"
$ws.Range("E43").Value = "  +0.89%  "

$ws.Range("D44").Value = "0.1666"
$ws.Range("E44").Value = "  +1.17%  "

$ws.Range("D45").Value = "8.584"
$ws.Range("E45").Value = "  +2.41%  "

$ws.Range("D46").Value = "0.4969"
$ws.Range("E46").Value = "  +2.32%  "

$ws.Range("D47").Value = "10.62"
$ws.Range("E47").Value = "  +2.04%  "

$ws.Range("D48").Value = "1.014"
$ws.Range("E48").Value = "  +1.10%  "

$ws.Range("D49").Value = "105.12"
$ws.Range("E49").Value = "  +2.53%  "

$ws.Range("D50").Value = "1.686"
$ws.Range("E50").Value = "  +1.63%  "

$ws.Range("D51").Value = "0.06328"
$ws.Range("E51").Value = "  +0.59%  "
